# Update capital structure database values for Kuwait Steel rows (2 and 3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Range("D$r").Value = -0.07099999999999999
    $ws.Range("E$r").ClearContents()

    $ws.Range("G$r").Value = 0.1142566191446029
    $ws.Range("H$r").Value = 0.1142566191446029
    $ws.Range("I$r").Value = 0.3529275524376054
    $ws.Range("J$r").Value = 0.3529275524376054
    $ws.Range("K$r").Value = -1.97
    $ws.Range("L$r").Value = -0.4012219959266802
    $ws.Range("M$r").Value = 15.93
    $ws.Range("N$r").Value = 0.2500784929356358
    $ws.Range("O$r").Value = -8.086294416243655
    $ws.Range("P$r").Value = 14.3
    $ws.Range("Q$r").Value = 0.2244897959183673
    $ws.Range("R$r").Value = -7.258883248730965
    $ws.Range("S$r").Value = 1.629999999999999
    $ws.Range("T$r").Value = 0.1023226616446955
    $ws.Range("U$r").Value = 3.39
    $ws.Range("V$r").Value = 0.0532182103610675
    $ws.Range("W$r").Value = -0.02739916550764951
    $ws.Range("X$r").Value = 0.05984447943376946
    $ws.Range("Y$r").Value = -0.08724364494141897
    $ws.Range("Z$r").Value = 0.09404215449650384
    $ws.Range("AA$r").Value = 0.03319006741241025
    $ws.Range("AB$r").Value = 0.0597486554290149
    $ws.Range("AC$r").Value = -0.02655858801660465

    $ws.Range("AE$r").Value = 0.2106285876567859
    $ws.Range("AF$r").Value = 0.2106285876567859
    $ws.Range("AG$r").Value = -3.179371412343214
    $ws.Range("AH$r").Value = 0.003295673854434052
    $ws.Range("AI$r").Value = 0.003391184285947543
    $ws.Range("AJ$r").Value = -0.05253368126767356
    $ws.Range("AK$r").Value = -0.05414402891817009
    $ws.Range("AL$r").Value = 0

    $ws.Range("AM$r").Value = -3.5

    $ws.Range("AO$r").ClearContents()
    $ws.Range("AP$r").Value = -1.704756789460169
    $ws.Range("AQ$r").Value = -0.4657142857142857
}
